# Apply the edits described by the diff:
#  - "RENAN NUNES" -> "RENAN" (3 occurrences: two body paragraphs + one table cell)
#  - "29" -> "18" and "abril" -> "Junho" in the signature date line
#     "Parnamirim/RN, 29 de abril de 2025."

$d = $word.ActiveDocument

# Replace every whole-word, case-matched occurrence of "RENAN NUNES" with "RENAN".
$range = $d.Content
$range.Find.Execute("RENAN NUNES", $true, $true, $false, $false, $false, `
                     $true, 1, $false, "RENAN", 2)

# Replace the day number in the signature date line. Use MatchWholeWord so we
# only hit the standalone "29" token (not e.g. the "2025" year or other IDs).
$range = $d.Content
$range.Find.Execute("29", $true, $true, $false, $false, $false, `
                     $true, 1, $false, "18", 2)

# Replace the month name in the signature date line.
$range = $d.Content
$range.Find.Execute("abril", $true, $true, $false, $false, $false, `
                     $true, 1, $false, "Junho", 2)
